# Corrige os dados e remove linhas de cabeçalho de subgrupo redundantes
# (situação do domicílio / grandes regiões), deslocando as linhas de
# dados para cima, e corrige o rótulo da coluna B da linha 2 para "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove a linha 5 ("situação do domicílio"), que ficava entre "brasil"
# e "urbana" sem dados proprios - as linhas abaixo sobem uma posicao.
$ws.Rows("5").Delete()

# Apos a remocao acima, a linha "grandes regioes" (antiga linha 8) passa
# a ser a linha 7 - remove-a tambem, deslocando "norte".."sul" para cima.
$ws.Rows("7").Delete()

# Corrige o rotulo da coluna "total" na segunda linha de cabecalho
# (antes continha o texto incorreto "unnamed: 1_level_1").
$ws.Range("B2").Value = "total"
